# Insert a new price-record row for "Albahaca" (Región Metropolitana) at
# row 351 of the "Vega Modelo de Temuco" sheet, pushing the existing rows
# 351-437 down to 352-438 (dimension grows from A1:R437 to A1:R438).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 351:437 down by one to make room for the new record.
$ws.Rows("351:351").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A351").Value = 10
$ws.Range("B351").Value = "Vega Modelo de Temuco"
$ws.Range("C351").Value = "La Araucanía"
$ws.Range("D351").Value = 45204
$ws.Range("E351").Value = 9
$ws.Range("F351").Value = 100112052
$ws.Range("G351").Value = "Albahaca"
$ws.Range("H351").Value = "Sin especificar"
$ws.Range("I351").Value = "Primera"
$ws.Range("J351").Value = 65
$ws.Range("K351").Value = 6000
$ws.Range("L351").Value = 6000
$ws.Range("M351").Value = 6000
$ws.Range("N351").Value = "$/paquete"
$ws.Range("O351").Value = "Región Metropolitana"
$ws.Range("P351").Value = 6000
$ws.Range("Q351").Value = 1
$ws.Range("R351").Value = "Hortaliza"
